$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-14 (columns B-G), with H as SUM(B:G) formula
$data = @(
    @(34, 228, 90, 122, 6, 7),    # row 2 - AGUADILLA
    @(17, 88, 51, 40, 2, 1),      # row 3 - AIBONITO
    @(34, 182, 97, 89, 3, 3),     # row 4 - ARECIBO
    @(58, 318, 179, 178, 9, 9),   # row 5 - BAYAMON
    @(45, 180, 108, 116, 6, 10),  # row 6 - CAGUAS
    @(30, 103, 65, 61, 5, 1),     # row 7 - CAROLINA
    @(29, 144, 55, 66, 8, 4),     # row 8 - FAJARDO
    @(24, 108, 58, 62, 7, 2),     # row 9 - GUAYAMA
    @(31, 255, 99, 100, 3, 2),    # row 10 - HUMACAO
    @(23, 60, 48, 55, 2, 2),      # row 11 - MAYAGUEZ
    @(46, 308, 136, 142, 2, 3),   # row 12 - PONCE
    @(67, 460, 223, 187, 15, 11), # row 13 - SAN JUAN
    @(12, 70, 26, 33, 1, 0)       # row 14 - UTUADO
)

$rowIndex = 2
foreach ($rowVals in $data) {
    $ws.Cells.Item($rowIndex, 2).Value = $rowVals[0]
    $ws.Cells.Item($rowIndex, 3).Value = $rowVals[1]
    $ws.Cells.Item($rowIndex, 4).Value = $rowVals[2]
    $ws.Cells.Item($rowIndex, 5).Value = $rowVals[3]
    $ws.Cells.Item($rowIndex, 6).Value = $rowVals[4]
    $ws.Cells.Item($rowIndex, 7).Value = $rowVals[5]
    $ws.Cells.Item($rowIndex, 8).Formula = "=SUM(B$rowIndex`:G$rowIndex)"
    $rowIndex++
}

# Update the active cell selection from B21 to G21
[void]$ws.Range("G21").Select()
